$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'22.388.80"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -3.98%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'1.570.28"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -3.50%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("E4").Value = "'  -0.19%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("E5").Value = "'  -0.10%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'289.15"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -2.82%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("E7").Value = "'  -2.05%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("E8").Value = "'  -2.21%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'0.3386"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -2.20%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'1.170"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -2.18%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'0.07626"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -4.74%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("E12").Value = "'  -0.13%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'21.24"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -2.42%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'6.061"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -3.37%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'6.909"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -3.90%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'1.572.09"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -3.62%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("D18").Value = "'89.51"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -5.13%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'0.06750"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -2.66%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("E20").Value = "'  -0.05%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'6.232"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -5.30%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'0.5335"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -5.70%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'16.53"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -3.92%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'11.99"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -2.40%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'22.412.82"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -3.94%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'2.364"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -2.99%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "'2.908"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -2.93%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("D28").Value = "'20.04"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -3.30%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").Value = "'145.41"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -3.71%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Value = "'4.963"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -3.73%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = "'125.44"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -4.38%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("D32").Value = "'1.748.77"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -3.65%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").Value = "'1.023"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +5.43%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").Value = "'6.239"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -6.42%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").Value = "'2.016"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -5.43%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("D36").Value = "'10.21"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -9.35%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").Value = "'0.08460"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -2.87%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("D38").Value = "'0.02536"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -4.27%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = "'0.2318"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -3.58%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("B40").Value = "'TrustWalletToken"
$ws.Range("B40").Style = "Normal"
$ws.Range("C40").Value = "'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("C40").Style = "Normal"
$ws.Range("D40").Value = "'1.337"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +3.43%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("B41").Value = "'InternetComputer(DFINITY)"
$ws.Range("B41").Style = "Normal"
$ws.Range("C41").Value = "'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("C41").Style = "Normal"
$ws.Range("D41").Value = "'5.529"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -4.85%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("B42").Value = "'Hedera"
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = "'0.06465"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -2.66%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'11.71"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -6.98%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'0.6352"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -5.85%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("E45").Value = "'  -7.93%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("E46").Value = "'  -0.05%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = "'0.5978"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -4.64%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'3.748"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -3.52%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = "'2.101"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -5.46%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = "'1.264"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +4.48%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = "'125.16"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -0.87%  "
$ws.Range("E51").Style = "Normal"
